# Fix a day/month transposition mistake in the "Timestamp" column (A) of
# every sheet: "03/02/2024 ..." (02 Mar, mis-typed as DD/MM) should read
# "02/03/2024 ..." (swap the day and month components, keep the time as-is).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $ts = $cell.Value2
        if ($ts -and $ts.ToString().StartsWith("03/02/2024")) {
            $cell.Value = $ts.ToString().Replace("03/02/2024", "02/03/2024")
        }
    }
}
